# Data update 2024-01-21
#
# Ticker-list maintenance on the FTSE 100 sheet:
#   - Dechra Pharmaceuticals (DPH) is removed (row 29); every row below it
#     shifts up by one to close the gap.
#   - Persimmon (PSN) is added back into the list in its correct
#     alphabetical slot, between "Pershing Square Holdings" and
#     "Phoenix Group" (ends up at row 68 after the shift above).
#
# Net effect: only rows 29-68 change; everything from row 69 down
# (Phoenix Group, Prudential, ...) keeps its original row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 currently holds DPH / Dechra Pharmaceuticals / Pharmaceuticals & Biotechnology.
# Deleting it shifts rows 30..101 up by one (row 30 DGE -> row 29, etc.).
$ws.Rows.Item(29).Delete()

# After the shift, "Pershing Square Holdings" sits at row 67 and
# "Phoenix Group" at row 68. Insert a fresh blank row at 68 to make room
# for Persimmon, pushing Phoenix Group (and everything after it) back
# down to its original row number.
$ws.Rows.Item(68).Insert()

$ws.Cells.Item(68, 1).Value = "PSN"
$ws.Cells.Item(68, 2).Value = "Persimmon"
$ws.Cells.Item(68, 3).Value = "Household Goods & Home Construction"
